# Auto-generated edit script: updates cryptocurrency price/volume data
# in Sheet1 to reflect the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 4).Value = '28.608.46'
$ws.Cells.Item(2, 5).Value = '  +1.64%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.828.03'
$ws.Cells.Item(3, 5).Value = '  +1.43%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.001'
$ws.Cells.Item(4, 5).Value = '  +0.07%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''316.10'
$ws.Cells.Item(5, 5).Value = '  -0.11%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.06%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.5344'
$ws.Cells.Item(7, 5).Value = '  -0.98%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.3989'
$ws.Cells.Item(8, 5).Value = '  +5.65%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.07778'
$ws.Cells.Item(9, 5).Value = '  +4.17%  '

# Row 10
$ws.Cells.Item(10, 2).Value = 'Polygon'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10, 4).Value = '''1.121'
$ws.Cells.Item(10, 5).Value = '  +2.46%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'OKB'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(11, 4).Value = '''41.96'
$ws.Cells.Item(11, 5).Value = '  -0.08%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''6.335'
$ws.Cells.Item(12, 5).Value = '  +2.05%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''21.12'
$ws.Cells.Item(13, 5).Value = '  +3.11%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''7.588'
$ws.Cells.Item(14, 5).Value = '  +2.85%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +0.07%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '1.830.11'
$ws.Cells.Item(16, 5).Value = '  +0.77%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''93.01'
$ws.Cells.Item(17, 5).Value = '  +3.63%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.00001091'
$ws.Cells.Item(18, 5).Value = '  +2.81%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''0.06561'
$ws.Cells.Item(19, 5).Value = '  +0.90%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''17.83'
$ws.Cells.Item(20, 5).Value = '  +2.73%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.04%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''6.090'
$ws.Cells.Item(22, 5).Value = '  +2.78%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '28.613.46'
$ws.Cells.Item(23, 5).Value = '  +1.54%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''11.22'
$ws.Cells.Item(24, 5).Value = '  +0.19%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''2.233'
$ws.Cells.Item(25, 5).Value = '  +6.75%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''20.85'
$ws.Cells.Item(26, 5).Value = '  +1.77%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''157.26'
$ws.Cells.Item(27, 5).Value = '  +0.81%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '2.049.81'
$ws.Cells.Item(28, 5).Value = '  +1.66%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''2.410'
$ws.Cells.Item(29, 5).Value = '  +3.87%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''125.65'
$ws.Cells.Item(30, 5).Value = '  +3.11%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''1.152'
$ws.Cells.Item(31, 5).Value = '  +2.99%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''0.1122'
$ws.Cells.Item(32, 5).Value = '  +0.86%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''5.747'
$ws.Cells.Item(33, 5).Value = '  +2.80%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.25%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.07342'
$ws.Cells.Item(35, 5).Value = '  +4.26%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''0.2274'
$ws.Cells.Item(36, 5).Value = '  +2.44%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''0.02357'
$ws.Cells.Item(37, 5).Value = '  +2.48%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''8.920'
$ws.Cells.Item(38, 5).Value = '  +5.52%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''5.219'
$ws.Cells.Item(39, 5).Value = '  +2.83%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''11.41'
$ws.Cells.Item(40, 5).Value = '  +2.88%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.6302'
$ws.Cells.Item(41, 5).Value = '  +2.23%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''1.198'
$ws.Cells.Item(42, 5).Value = '  +1.87%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -2.99%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''13.59'
$ws.Cells.Item(45, 5).Value = '  +1.62%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''0.5943'
$ws.Cells.Item(46, 5).Value = '  +3.24%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''3.711'
$ws.Cells.Item(47, 5).Value = '  +0.80%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''125.98'
$ws.Cells.Item(48, 5).Value = '  +0.68%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''2.000'
$ws.Cells.Item(49, 5).Value = '  +3.97%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''1.195'
$ws.Cells.Item(50, 5).Value = '  +0.56%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''0.06954'
$ws.Cells.Item(51, 5).Value = '  +1.99%  '

